$d = $word.ActiveDocument

$d.Content.Find.Execute("679÷7=97, 0", $true, $false, $false, $false, $false, $true, 1, $false, "979÷7=139, 6", 2) | Out-Null
$d.Content.Find.Execute("668÷3=222, 2", $true, $false, $false, $false, $false, $true, 1, $false, "346÷9=38, 4", 2) | Out-Null
$d.Content.Find.Execute("812÷8=101, 4", $true, $false, $false, $false, $false, $true, 1, $false, "293÷6=48, 5", 2) | Out-Null
$d.Content.Find.Execute("215÷4=53, 3", $true, $false, $false, $false, $false, $true, 1, $false, "219÷6=36, 3", 2) | Out-Null
$d.Content.Find.Execute("636÷5=127, 1", $true, $false, $false, $false, $false, $true, 1, $false, "615÷4=153, 3", 2) | Out-Null
$d.Content.Find.Execute("758÷7=108, 2", $true, $false, $false, $false, $false, $true, 1, $false, "746÷4=186, 2", 2) | Out-Null
$d.Content.Find.Execute("684÷4=171, 0", $true, $false, $false, $false, $false, $true, 1, $false, "664÷9=73, 7", 2) | Out-Null
$d.Content.Find.Execute("529÷8=66, 1", $true, $false, $false, $false, $false, $true, 1, $false, "923÷7=131, 6", 2) | Out-Null
$d.Content.Find.Execute("299÷7=42, 5", $true, $false, $false, $false, $false, $true, 1, $false, "615÷3=205, 0", 2) | Out-Null
$d.Content.Find.Execute("933÷8=116, 5", $true, $false, $false, $false, $false, $true, 1, $false, "276÷4=69, 0", 2) | Out-Null
$d.Content.Find.Execute("824÷3=274, 2", $true, $false, $false, $false, $false, $true, 1, $false, "403÷4=100, 3", 2) | Out-Null
$d.Content.Find.Execute("851÷4=212, 3", $true, $false, $false, $false, $false, $true, 1, $false, "695÷3=231, 2", 2) | Out-Null
$d.Content.Find.Execute("975÷7=139, 2", $true, $false, $false, $false, $false, $true, 1, $false, "637÷2=318, 1", 2) | Out-Null
$d.Content.Find.Execute("146÷2=73, 0", $true, $false, $false, $false, $false, $true, 1, $false, "766÷5=153, 1", 2) | Out-Null
$d.Content.Find.Execute("465÷5=93, 0", $true, $false, $false, $false, $false, $true, 1, $false, "153÷3=51, 0", 2) | Out-Null
$d.Content.Find.Execute("499÷4=124, 3", $true, $false, $false, $false, $false, $true, 1, $false, "549÷6=91, 3", 2) | Out-Null
$d.Content.Find.Execute("497÷9=55, 2", $true, $false, $false, $false, $false, $true, 1, $false, "651÷2=325, 1", 2) | Out-Null
$d.Content.Find.Execute("892÷3=297, 1", $true, $false, $false, $false, $false, $true, 1, $false, "104÷2=52, 0", 2) | Out-Null
$d.Content.Find.Execute("538÷4=134, 2", $true, $false, $false, $false, $false, $true, 1, $false, "496÷8=62, 0", 2) | Out-Null
$d.Content.Find.Execute("619÷9=68, 7", $true, $false, $false, $false, $false, $true, 1, $false, "633÷7=90, 3", 2) | Out-Null
$d.Content.Find.Execute("499÷3=166, 1", $true, $false, $false, $false, $false, $true, 1, $false, "594÷8=74, 2", 2) | Out-Null
$d.Content.Find.Execute("353÷5=70, 3", $true, $false, $false, $false, $false, $true, 1, $false, "376÷8=47, 0", 2) | Out-Null
$d.Content.Find.Execute("455÷6=75, 5", $true, $false, $false, $false, $false, $true, 1, $false, "901÷8=112, 5", 2) | Out-Null
$d.Content.Find.Execute("502÷7=71, 5", $true, $false, $false, $false, $false, $true, 1, $false, "262÷8=32, 6", 2) | Out-Null
$d.Content.Find.Execute("522÷6=87, 0", $true, $false, $false, $false, $false, $true, 1, $false, "264÷6=44, 0", 2) | Out-Null
